$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '68.748.45'
$cell.ClearFormats()
$ws.Range("E2").Value = '  +0.87%  '

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '2.733.38'
$cell.ClearFormats()
$ws.Range("E3").Value = '  +3.61%  '

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.ClearFormats()
$ws.Range("E4").Value = '  -0.03%  '

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '602.33'
$cell.ClearFormats()
$ws.Range("E5").Value = '  +1.35%  '

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '168.92'
$cell.ClearFormats()
$ws.Range("E6").Value = '  +6.68%  '

$ws.Range("E7").Value = '  -0.01%  '

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '0.546'
$cell.ClearFormats()
$ws.Range("E8").Value = '  +0.64%  '

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '2.732.73'
$cell.ClearFormats()
$ws.Range("E9").Value = '  +3.61%  '

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '0.144'
$cell.ClearFormats()
$ws.Range("E10").Value = '  +1.93%  '

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '0.366'
$cell.ClearFormats()
$ws.Range("E11").Value = '  +4.86%  '

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '5.33'
$cell.ClearFormats()
$ws.Range("E12").Value = '  +1.39%  '

$ws.Range("E13").Value = '  -0.25%  '

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '28.76'
$cell.ClearFormats()
$ws.Range("E14").Value = '  +3.27%  '

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '3.235.32'
$cell.ClearFormats()
$ws.Range("E15").Value = '  +3.71%  '

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '0.0000190'
$cell.ClearFormats()
$ws.Range("E16").Value = '  +1.44%  '

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '68.682.87'
$cell.ClearFormats()
$ws.Range("E17").Value = '  +0.96%  '

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '2.730.63'
$cell.ClearFormats()
$ws.Range("E18").Value = '  +2.87%  '

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '11.83'
$cell.ClearFormats()
$ws.Range("E19").Value = '  +4.44%  '

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '371.07'
$cell.ClearFormats()
$ws.Range("E20").Value = '  +3.53%  '

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '7.71'
$cell.ClearFormats()
$ws.Range("E21").Value = '  +5.60%  '

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '4.53'
$cell.ClearFormats()
$ws.Range("E22").Value = '  +3.67%  '

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '4.97'
$cell.ClearFormats()
$ws.Range("E23").Value = '  +4.63%  '

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '2.13'
$cell.ClearFormats()
$ws.Range("E24").Value = '  +3.23%  '

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '73.93'
$cell.ClearFormats()
$ws.Range("E25").Value = '  -0.99%  '

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.ClearFormats()
$ws.Range("E26").Value = '  +0.02%  '

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '9.94'
$cell.ClearFormats()
$ws.Range("E27").Value = '  +1.94%  '

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '2.876.84'
$cell.ClearFormats()
$ws.Range("E28").Value = '  +3.64%  '

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '0.0000106'
$cell.ClearFormats()
$ws.Range("E29").Value = '  +3.74%  '

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '596.38'
$cell.ClearFormats()
$ws.Range("E30").Value = '  +6.36%  '

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '0.998'
$cell.ClearFormats()
$ws.Range("E31").Value = '  -0.12%  '

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '8.30'
$cell.ClearFormats()
$ws.Range("E32").Value = '  +4.71%  '

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '1.45'
$cell.ClearFormats()

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '1.96'
$cell.ClearFormats()
$ws.Range("E34").Value = '  +6.59%  '

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '0.133'
$cell.ClearFormats()
$ws.Range("E35").Value = '  +4.63%  '

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '1.63'
$cell.ClearFormats()
$ws.Range("E36").Value = '  +5.67%  '

$ws.Range("E37").Value = '  -0.08%  '

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '162.29'
$cell.ClearFormats()
$ws.Range("E38").Value = '  +1.65%  '

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '19.82'
$cell.ClearFormats()
$ws.Range("E39").Value = '  +0.69%  '

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '0.382'
$cell.ClearFormats()
$ws.Range("E40").Value = '  +3.56%  '

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '1.92'
$cell.ClearFormats()
$ws.Range("E41").Value = '  +2.93%  '

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '5.50'
$cell.ClearFormats()
$ws.Range("E42").Value = '  +4.52%  '

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '2.69'
$cell.ClearFormats()
$ws.Range("E43").Value = '  +4.17%  '

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '17.98'
$cell.ClearFormats()
$ws.Range("E44").Value = '  +1.07%  '

$ws.Range("E45").Value = '  -0.01%  '

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '0.0₆0315'
$cell.ClearFormats()
$ws.Range("E46").Value = '  -3.04%  '

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '40.91'
$cell.ClearFormats()
$ws.Range("E47").Value = '  +1.73%  '

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '156.83'
$cell.ClearFormats()
$ws.Range("E48").Value = '  +0.10%  '

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '3.97'
$cell.ClearFormats()
$ws.Range("E49").Value = '  +6.72%  '

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '1.81'
$cell.ClearFormats()
$ws.Range("E50").Value = '  +7.82%  '

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '0.605'
$cell.ClearFormats()
$ws.Range("E51").Value = '  +7.63%  '
